# edit.ps1 - applies the commit's changes to the document via Word COM-interop
$d = $word.ActiveDocument

# Helper: split a run in two at an exact character offset by briefly adding
# (then immediately removing) a bookmark there. Word always breaks a run at
# a bookmark boundary, and deleting the bookmark afterwards leaves the break
# in place without any other residue.
function Split-AtOffset($absOffset, $markName) {
    $r = $d.Range($absOffset, $absOffset)
    $d.Bookmarks.Add($markName, $r) | Out-Null
    $d.Bookmarks($markName).Delete()
}

# Helper: find searchText (first match in the document) and split right
# after the first $offsetFromStart characters of the match.
function Split-AfterMatch($searchText, $offsetFromStart, $markName) {
    $rng = $d.Content
    $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0) | Out-Null
    if ($rng.Find.Found) {
        Split-AtOffset ($rng.Start + $offsetFromStart) $markName
    }
}

# ---------------------------------------------------------------------
# Hunk 1: "...First, we analyze the top 10 movies genre on each year..."
#   -> drop " genre" and split the sentence into three runs.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("top 10 movies genre on each year", $true, $false, $false, `
                   $false, $false, $true, 1, $false, `
                   "top 10 movies on each year", 2) | Out-Null

# The text edit above causes Word to coalesce every same-formatted run in
# the paragraph, so re-establish the two boundaries that already existed
# (before "derive the conclusion" and before the closing ".") as well as
# the two new boundaries the commit introduces.
Split-AfterMatch "Third, derive the conclusion" ("Third, ").Length "TmpSplitA"
Split-AfterMatch "genres of the movies." ("genres of the movies").Length "TmpSplitB"
Split-AfterMatch "First, we analyze the top 10 movies " ("First, we ").Length "TmpSplitC"
Split-AfterMatch "analyze the top 10 movies on each year" `
    ("analyze the top 10 movies ").Length "TmpSplitD"

# ---------------------------------------------------------------------
# Hunk 2: "In the application layer, we will use Python..."
#   -> "we used" (tense change) and split into three runs.
# ---------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("In the application layer, we will use Python", $true, `
                    $false, $false, $false, $false, $true, 1, $false, `
                    "In the application layer, we used Python", 2) | Out-Null

Split-AfterMatch "In the application layer, we used Python" `
    ("In the application layer, we use").Length "TmpSplitE"
Split-AfterMatch "we used Python" ("we use" + "d").Length "TmpSplitF"

# ---------------------------------------------------------------------
# Hunk 3 + 4: "top 10 movie genre on each year (sophisticated query)"
#   -> split into two runs with a "_GoBack" bookmark in between.
# "_GoBack" is Word's special single-instance auto bookmark: (re)adding it
# here automatically relocates it away from its old spot just before
# "Fig. Trend for genre ..." (hunk 4), removing the old bookmarkStart/End
# pair with no extra work needed.
# ---------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("top 10 movie genre", $true, $false, $false, $false, `
                    $false, $true, 1, $false, "", 0) | Out-Null
if ($rng3.Find.Found) {
    $splitPoint = $rng3.End
    $bmRange = $d.Range($splitPoint, $splitPoint)
    $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
}
